# Refresh "Top Gainers" leaderboard (re-sorted by Latest % desc) and fix one cell on "Top Losers".
$wb = $excel.ActiveWorkbook
$wsGainers = $wb.Worksheets.Item("Top Gainers")
$wsLosers = $wb.Worksheets.Item("Top Losers")

$wsGainers.Cells.Item(2, 2).Value = "BLUEDART"
$wsGainers.Cells.Item(2, 3).Value = 14.2935
$wsGainers.Cells.Item(2, 4).Value = 13.7184
$wsGainers.Cells.Item(2, 5).Value = 10.9767
$wsGainers.Cells.Item(3, 2).Value = "SOLARWORLD"
$wsGainers.Cells.Item(3, 3).Value = 14.238
$wsGainers.Cells.Item(3, 4).Value = 10.2449
$wsGainers.Cells.Item(3, 5).Value = 5.7805
$wsGainers.Cells.Item(4, 2).Value = "DREDGECORP"
$wsGainers.Cells.Item(4, 3).Value = 12.5116
$wsGainers.Cells.Item(4, 4).Value = 17.0315
$wsGainers.Cells.Item(4, 5).Value = 17.7966
$wsGainers.Cells.Item(5, 2).Value = "ADANIGREEN"
$wsGainers.Cells.Item(5, 3).Value = 11.8303
$wsGainers.Cells.Item(5, 4).Value = 8.9128
$wsGainers.Cells.Item(5, 5).Value = 9.3636
$wsGainers.Cells.Item(7, 2).Value = "VBL"
$wsGainers.Cells.Item(7, 3).Value = 9.4132
$wsGainers.Cells.Item(7, 4).Value = 7.694
$wsGainers.Cells.Item(7, 5).Value = 11.9901
$wsGainers.Cells.Item(8, 2).Value = "HEG"
$wsGainers.Cells.Item(8, 3).Value = 8.1713
$wsGainers.Cells.Item(8, 4).Value = 12.1776
$wsGainers.Cells.Item(8, 5).Value = 14.5771
$wsGainers.Cells.Item(9, 2).Value = "BUTTERFLY"
$wsGainers.Cells.Item(9, 3).Value = 7.6763
$wsGainers.Cells.Item(9, 4).Value = 10.5879
$wsGainers.Cells.Item(9, 5).Value = 13.1883
$wsGainers.Cells.Item(10, 2).Value = "WALCHANNAG"
$wsGainers.Cells.Item(10, 3).Value = 7.1578
$wsGainers.Cells.Item(10, 4).Value = 4.683
$wsGainers.Cells.Item(10, 5).Value = -4.5619
$wsGainers.Cells.Item(12, 2).Value = "ABREL"
$wsGainers.Cells.Item(12, 3).Value = 6.7686
$wsGainers.Cells.Item(12, 4).Value = 7.5428
$wsGainers.Cells.Item(12, 5).Value = 7.0972
$wsGainers.Cells.Item(13, 2).Value = "SAIL"
$wsGainers.Cells.Item(13, 3).Value = 6.5602
$wsGainers.Cells.Item(13, 4).Value = 8.7826
$wsGainers.Cells.Item(13, 5).Value = 4.7219
$wsGainers.Cells.Item(14, 2).Value = "SANDUMA"
$wsGainers.Cells.Item(14, 3).Value = 6.4045
$wsGainers.Cells.Item(14, 4).Value = 3.9096
$wsGainers.Cells.Item(14, 5).Value = 32.5377
$wsGainers.Cells.Item(15, 2).Value = "JISLJALEQS"
$wsGainers.Cells.Item(15, 3).Value = 6.2627
$wsGainers.Cells.Item(15, 4).Value = 5.5544
$wsGainers.Cells.Item(15, 5).Value = -0.5519
$wsGainers.Cells.Item(16, 2).Value = "M&MFIN"
$wsGainers.Cells.Item(16, 3).Value = 6.2021
$wsGainers.Cells.Item(16, 4).Value = 6.8254
$wsGainers.Cells.Item(16, 5).Value = 15.692
$wsGainers.Cells.Item(17, 2).Value = "FISCHER"
$wsGainers.Cells.Item(17, 3).Value = 6.0873
$wsGainers.Cells.Item(17, 4).Value = 11.0876
$wsGainers.Cells.Item(17, 5).Value = 4.2523
$wsGainers.Cells.Item(18, 2).Value = "ADANIENSOL"
$wsGainers.Cells.Item(18, 3).Value = 6.0465
$wsGainers.Cells.Item(18, 4).Value = 3.4578
$wsGainers.Cells.Item(18, 5).Value = 11.9977
$wsGainers.Cells.Item(19, 2).Value = "UTKARSHBNK"
$wsGainers.Cells.Item(19, 3).Value = 5.8522
$wsGainers.Cells.Item(19, 4).Value = -5.0207
$wsGainers.Cells.Item(19, 5).Value = -1.7159
$wsGainers.Cells.Item(20, 2).Value = "POKARNA"
$wsGainers.Cells.Item(20, 3).Value = 5.8251
$wsGainers.Cells.Item(20, 4).Value = -1.5953
$wsGainers.Cells.Item(20, 5).Value = 18.6047
$wsGainers.Cells.Item(21, 2).Value = "ABDL"
$wsGainers.Cells.Item(21, 3).Value = 5.6218
$wsGainers.Cells.Item(21, 4).Value = 4.5201
$wsGainers.Cells.Item(21, 5).Value = 27.2996
$wsGainers.Cells.Item(22, 2).Value = "CELLO"
$wsGainers.Cells.Item(22, 3).Value = 5.6032
$wsGainers.Cells.Item(22, 4).Value = 4.4395
$wsGainers.Cells.Item(22, 5).Value = 14.3576
$wsGainers.Cells.Item(23, 2).Value = "EPACKPEB"
$wsGainers.Cells.Item(23, 3).Value = 5.5999
$wsGainers.Cells.Item(23, 4).Value = -0.9271
$wsGainers.Cells.Item(23, 5).Value = "N/A"
$wsGainers.Cells.Item(24, 2).Value = "VAIBHAVGBL"
$wsGainers.Cells.Item(24, 3).Value = 5.556
$wsGainers.Cells.Item(24, 4).Value = 6.0822
$wsGainers.Cells.Item(24, 5).Value = 12.4107
$wsGainers.Cells.Item(25, 2).Value = "GRAPHITE"
$wsGainers.Cells.Item(25, 3).Value = 5.5311
$wsGainers.Cells.Item(25, 4).Value = 11.9508
$wsGainers.Cells.Item(25, 5).Value = 12.1523
$wsGainers.Cells.Item(26, 2).Value = "FIVESTAR"
$wsGainers.Cells.Item(26, 3).Value = 5.2509
$wsGainers.Cells.Item(26, 4).Value = 5.2607
$wsGainers.Cells.Item(26, 5).Value = 5.3392
$wsGainers.Cells.Item(27, 2).Value = "IOC"
$wsGainers.Cells.Item(27, 3).Value = 5.1255
$wsGainers.Cells.Item(27, 4).Value = 8.0269
$wsGainers.Cells.Item(27, 5).Value = 8.4452
$wsGainers.Cells.Item(28, 2).Value = "SURYAROSNI"
$wsGainers.Cells.Item(28, 3).Value = 5.0721
$wsGainers.Cells.Item(28, 4).Value = 11.4804
$wsGainers.Cells.Item(28, 5).Value = 3.1086
$wsGainers.Cells.Item(29, 2).Value = "MEGASOFT"
$wsGainers.Cells.Item(29, 3).Value = 4.9974
$wsGainers.Cells.Item(29, 4).Value = 15.7588
$wsGainers.Cells.Item(29, 5).Value = 33.5271
$wsGainers.Cells.Item(30, 2).Value = "PROZONER"
$wsGainers.Cells.Item(30, 3).Value = 4.9921
$wsGainers.Cells.Item(30, 4).Value = 15.7468
$wsGainers.Cells.Item(30, 5).Value = 36.095
$wsGainers.Cells.Item(31, 2).Value = "ATGL"
$wsGainers.Cells.Item(31, 3).Value = 4.9919
$wsGainers.Cells.Item(31, 4).Value = 4.7305
$wsGainers.Cells.Item(31, 5).Value = 4.27
$wsGainers.Cells.Item(32, 2).Value = "STALLION"
$wsGainers.Cells.Item(32, 3).Value = 4.9914
$wsGainers.Cells.Item(32, 4).Value = -5.2229
$wsGainers.Cells.Item(32, 5).Value = 21.4391
$wsGainers.Cells.Item(33, 2).Value = "INDOTHAI"
$wsGainers.Cells.Item(33, 3).Value = 4.9883
$wsGainers.Cells.Item(33, 4).Value = 4.7163
$wsGainers.Cells.Item(33, 5).Value = 43.9974
$wsGainers.Cells.Item(34, 2).Value = "HITECHGEAR"
$wsGainers.Cells.Item(34, 3).Value = 4.8651
$wsGainers.Cells.Item(34, 4).Value = 2.1287
$wsGainers.Cells.Item(34, 5).Value = 10.9905
$wsGainers.Cells.Item(35, 2).Value = "BAJAJINDEF"
$wsGainers.Cells.Item(35, 3).Value = 4.7737
$wsGainers.Cells.Item(35, 4).Value = 3.6917
$wsGainers.Cells.Item(35, 5).Value = 10.7236
$wsGainers.Cells.Item(36, 2).Value = "DATAMATICS"
$wsGainers.Cells.Item(36, 3).Value = 4.63
$wsGainers.Cells.Item(36, 4).Value = 7.0385
$wsGainers.Cells.Item(36, 5).Value = 15.4314
$wsGainers.Cells.Item(37, 2).Value = "CMSINFO"
$wsGainers.Cells.Item(37, 3).Value = 4.5725
$wsGainers.Cells.Item(37, 4).Value = 3.3422
$wsGainers.Cells.Item(37, 5).Value = 3.5499
$wsGainers.Cells.Item(38, 2).Value = "JKIL"
$wsGainers.Cells.Item(38, 3).Value = 4.3862
$wsGainers.Cells.Item(38, 4).Value = 3.1925
$wsGainers.Cells.Item(38, 5).Value = 2.0017
$wsGainers.Cells.Item(39, 2).Value = "GMBREW"
$wsGainers.Cells.Item(39, 3).Value = 4.3543
$wsGainers.Cells.Item(39, 4).Value = -0.0949
$wsGainers.Cells.Item(39, 5).Value = 79.8121
$wsGainers.Cells.Item(40, 2).Value = "AXISCADES"
$wsGainers.Cells.Item(40, 3).Value = 4.3405
$wsGainers.Cells.Item(40, 4).Value = 6.8028
$wsGainers.Cells.Item(40, 5).Value = -3.1807
$wsGainers.Cells.Item(41, 2).Value = "GENUSPOWER"
$wsGainers.Cells.Item(41, 3).Value = 4.3243
$wsGainers.Cells.Item(41, 4).Value = 2.6425
$wsGainers.Cells.Item(41, 5).Value = -0.4032
$wsGainers.Cells.Item(42, 2).Value = "SAMBHV"
$wsGainers.Cells.Item(42, 3).Value = 4.3166
$wsGainers.Cells.Item(42, 4).Value = 2.8031
$wsGainers.Cells.Item(42, 5).Value = 5.3506
$wsGainers.Cells.Item(43, 2).Value = "STLTECH"
$wsGainers.Cells.Item(43, 3).Value = 4.2953
$wsGainers.Cells.Item(43, 4).Value = 1.299
$wsGainers.Cells.Item(43, 5).Value = 7.3983
$wsGainers.Cells.Item(44, 2).Value = "SGMART"
$wsGainers.Cells.Item(44, 3).Value = 4.2877
$wsGainers.Cells.Item(44, 4).Value = 8.2735
$wsGainers.Cells.Item(44, 5).Value = 2.552
$wsGainers.Cells.Item(45, 2).Value = "TMB"
$wsGainers.Cells.Item(45, 3).Value = 4.2788
$wsGainers.Cells.Item(45, 4).Value = 8.0084
$wsGainers.Cells.Item(45, 5).Value = 15.2269
$wsGainers.Cells.Item(46, 2).Value = "GPIL"
$wsGainers.Cells.Item(46, 3).Value = 4.2756
$wsGainers.Cells.Item(46, 4).Value = 6.4483
$wsGainers.Cells.Item(46, 5).Value = 14.5711
$wsGainers.Cells.Item(47, 2).Value = "PROSTARM"
$wsGainers.Cells.Item(47, 3).Value = 4.2704
$wsGainers.Cells.Item(47, 4).Value = 1.3783
$wsGainers.Cells.Item(47, 5).Value = -7.5853
$wsGainers.Cells.Item(48, 2).Value = "SUNFLAG"
$wsGainers.Cells.Item(48, 3).Value = 4.2432
$wsGainers.Cells.Item(48, 4).Value = 4.58
$wsGainers.Cells.Item(48, 5).Value = 4.879
$wsGainers.Cells.Item(49, 2).Value = "SRM"
$wsGainers.Cells.Item(49, 3).Value = 4.1267
$wsGainers.Cells.Item(49, 4).Value = 3.8172
$wsGainers.Cells.Item(49, 5).Value = 4.7316
$wsGainers.Cells.Item(50, 2).Value = "PDSL"
$wsGainers.Cells.Item(50, 3).Value = 4.0782
$wsGainers.Cells.Item(50, 4).Value = 2.1122
$wsGainers.Cells.Item(50, 5).Value = 7.8894
$wsGainers.Cells.Item(51, 2).Value = "TCI"
$wsGainers.Cells.Item(51, 3).Value = 3.9632
$wsGainers.Cells.Item(51, 4).Value = 3.8669
$wsGainers.Cells.Item(51, 5).Value = 4.3681
$wsGainers.Cells.Item(52, 2).Value = "MRPL"
$wsGainers.Cells.Item(52, 3).Value = 3.9568
$wsGainers.Cells.Item(52, 4).Value = 9.3868
$wsGainers.Cells.Item(52, 5).Value = 19.7003
$wsGainers.Cells.Item(53, 2).Value = "GPPL"
$wsGainers.Cells.Item(53, 3).Value = 3.9562
$wsGainers.Cells.Item(53, 4).Value = 2.9525
$wsGainers.Cells.Item(53, 5).Value = 4.5877
$wsGainers.Cells.Item(54, 2).Value = "LLOYDSENT"
$wsGainers.Cells.Item(54, 3).Value = 3.9018
$wsGainers.Cells.Item(54, 4).Value = 1.1884
$wsGainers.Cells.Item(54, 5).Value = 10.5288
$wsGainers.Cells.Item(55, 2).Value = "VINCOFE"
$wsGainers.Cells.Item(55, 3).Value = 3.8925
$wsGainers.Cells.Item(55, 4).Value = 10.7668
$wsGainers.Cells.Item(55, 5).Value = 9.1395
$wsGainers.Cells.Item(56, 2).Value = "PRAKASH"
$wsGainers.Cells.Item(56, 3).Value = 3.8648
$wsGainers.Cells.Item(56, 4).Value = 4.7712
$wsGainers.Cells.Item(56, 5).Value = 1.5084
$wsGainers.Cells.Item(57, 2).Value = "ICRA"
$wsGainers.Cells.Item(57, 3).Value = 3.862
$wsGainers.Cells.Item(57, 4).Value = 4.5433
$wsGainers.Cells.Item(57, 5).Value = 2.9458
$wsGainers.Cells.Item(58, 2).Value = "NBCC"
$wsGainers.Cells.Item(58, 3).Value = 3.8618
$wsGainers.Cells.Item(58, 4).Value = 2.5786
$wsGainers.Cells.Item(58, 5).Value = 6.9948
$wsGainers.Cells.Item(59, 2).Value = "RECLTD"
$wsGainers.Cells.Item(59, 3).Value = 3.8352
$wsGainers.Cells.Item(59, 4).Value = 2.818
$wsGainers.Cells.Item(59, 5).Value = 2.7491
$wsGainers.Cells.Item(60, 2).Value = "HCC"
$wsGainers.Cells.Item(60, 3).Value = 3.7446
$wsGainers.Cells.Item(60, 4).Value = 2.6464
$wsGainers.Cells.Item(60, 5).Value = 7.3828
$wsGainers.Cells.Item(61, 2).Value = "AVALON"
$wsGainers.Cells.Item(61, 3).Value = 3.7195
$wsGainers.Cells.Item(61, 4).Value = 8.4955
$wsGainers.Cells.Item(61, 5).Value = 20.4776
$wsGainers.Cells.Item(62, 2).Value = "INDORAMA"
$wsGainers.Cells.Item(62, 3).Value = 3.687
$wsGainers.Cells.Item(62, 4).Value = 2.8512
$wsGainers.Cells.Item(62, 5).Value = 13.9974
$wsGainers.Cells.Item(63, 2).Value = "STAR"
$wsGainers.Cells.Item(63, 3).Value = 3.6855
$wsGainers.Cells.Item(63, 4).Value = 3.6155
$wsGainers.Cells.Item(63, 5).Value = 2.8516
$wsGainers.Cells.Item(64, 2).Value = "RAJRATAN"
$wsGainers.Cells.Item(64, 3).Value = 3.6659
$wsGainers.Cells.Item(64, 4).Value = 1.2098
$wsGainers.Cells.Item(64, 5).Value = 27.2943
$wsGainers.Cells.Item(65, 2).Value = "GAIL"
$wsGainers.Cells.Item(65, 3).Value = 3.6591
$wsGainers.Cells.Item(65, 4).Value = 2.1931
$wsGainers.Cells.Item(65, 5).Value = 4.9351
$wsGainers.Cells.Item(66, 2).Value = "DCW"
$wsGainers.Cells.Item(66, 3).Value = 3.5953
$wsGainers.Cells.Item(66, 4).Value = 2.165
$wsGainers.Cells.Item(66, 5).Value = -4.1225
$wsGainers.Cells.Item(67, 2).Value = "ASHAPURMIN"
$wsGainers.Cells.Item(67, 3).Value = 3.587
$wsGainers.Cells.Item(67, 4).Value = 6.238
$wsGainers.Cells.Item(67, 5).Value = 2.0277
$wsGainers.Cells.Item(68, 2).Value = "MAITHANALL"
$wsGainers.Cells.Item(68, 3).Value = 3.579
$wsGainers.Cells.Item(68, 4).Value = 2.6205
$wsGainers.Cells.Item(68, 5).Value = 1.8923
$wsGainers.Cells.Item(69, 2).Value = "DCMSHRIRAM"
$wsGainers.Cells.Item(69, 3).Value = 3.5092
$wsGainers.Cells.Item(69, 4).Value = 10.0678
$wsGainers.Cells.Item(69, 5).Value = 17.4318
$wsGainers.Cells.Item(70, 2).Value = "MSPL"
$wsGainers.Cells.Item(70, 3).Value = 3.4967
$wsGainers.Cells.Item(70, 4).Value = 2.1534
$wsGainers.Cells.Item(70, 5).Value = -5.1233
$wsGainers.Cells.Item(72, 2).Value = "ORIENTTECH"
$wsGainers.Cells.Item(72, 3).Value = 3.4509
$wsGainers.Cells.Item(72, 4).Value = 0.1606
$wsGainers.Cells.Item(72, 5).Value = 32.1979
$wsGainers.Cells.Item(73, 2).Value = "HLEGLAS"
$wsGainers.Cells.Item(73, 3).Value = 3.432
$wsGainers.Cells.Item(73, 4).Value = 7.8788
$wsGainers.Cells.Item(73, 5).Value = 26.8456
$wsGainers.Cells.Item(74, 2).Value = "HARSHA"
$wsGainers.Cells.Item(74, 3).Value = 3.3798
$wsGainers.Cells.Item(74, 4).Value = 2.3129
$wsGainers.Cells.Item(74, 5).Value = 5.016
$wsGainers.Cells.Item(75, 2).Value = "RESPONIND"
$wsGainers.Cells.Item(75, 3).Value = 3.3742
$wsGainers.Cells.Item(75, 4).Value = 7.3779
$wsGainers.Cells.Item(75, 5).Value = 6.2052
$wsGainers.Cells.Item(76, 2).Value = "AWHCL"
$wsGainers.Cells.Item(76, 3).Value = 3.3511
$wsGainers.Cells.Item(76, 4).Value = 3.5878
$wsGainers.Cells.Item(76, 5).Value = 0.4534

$wsLosers.Cells.Item(12, 4).Value = 5.9123
